$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
Write-Host "Before insert J1:" $ws.Cells.Item(1,10).Value
$ws.Columns.Item(10).Insert()
Write-Host "After insert J1 (should be empty):" $ws.Cells.Item(1,10).Value
Write-Host "After insert K1 (should be old J1=7):" $ws.Cells.Item(1,11).Value
